$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------------
# Constants sheet: rework the rows to add the new queue / retry / exception
# related settings (MaxConsecutiveSystemExceptions, RetryNumberGetTransactionItem,
# RetryNumberSetTransactionStatus, ShouldMarkJobAsFaulted, ExceptionMessage_
# ConsecutiveErrors) and move MailException up under the new header row.
# ---------------------------------------------------------------------------

# Row 2 becomes a plain header row (Name / Value / Description)
$wsConstants.Range("A2").Value = "Name"
$wsConstants.Range("B2").Value = "Value"
$wsConstants.Range("C2").Value = "Description"

# Row 3: MailException (moved up from its old row 12 spot, e-mail corrected)
$wsConstants.Range("A3").Value = "MailException"
$wsConstants.Range("B3").Value = "parnupong.k@thaibev.com"
$wsConstants.Range("C3").ClearContents()

# Row 4: MaxRetryNumber
$wsConstants.Range("A4").Value = "MaxRetryNumber"
$wsConstants.Range("B4").Value = 0
$wsConstants.Range("C4").Value = "Must be 0 if working with Orchestrator queues. If > 0, the robot will retry the same transaction which failed with a system exception. Must be an integer value."

# Row 5: MaxConsecutiveSystemExceptions (new)
$wsConstants.Range("A5").Value = "MaxConsecutiveSystemExceptions"
$wsConstants.Range("B5").Value = 0
$wsConstants.Range("C5").Value = "The number of consecutive system exceptions allowed. If MaxConsecutiveSystemExceptions is reached, the job is stopped. To disable this feature, set the value to 0. "

# Row 6: blank separator (was LogMessage_GetTransactionData before, now empty)
$wsConstants.Range("A6:C6").ClearContents()

# Row 7: ExScreenshotsFolderPath
$wsConstants.Range("A7").Value = "ExScreenshotsFolderPath"
$wsConstants.Range("B7").Value = "Exceptions_Screenshots"
$wsConstants.Range("C7").Value = "Where to save exceptions screenshots - can be a full or a relative path."

# Row 8: blank separator (was LogMessage_Success before, now empty)
$wsConstants.Range("A8:C8").ClearContents()

# Row 9: LogMessage_GetTransactionData
$wsConstants.Range("A9").Value = "LogMessage_GetTransactionData"
$wsConstants.Range("B9").Value = "Processing Transaction Number: "
$wsConstants.Range("C9").Value = "Static part of logging message. Calling Get Transaction Data."

# Row 10: LogMessage_GetTransactionDataError
$wsConstants.Range("A10").Value = "LogMessage_GetTransactionDataError"
$wsConstants.Range("B10").Value = "Error getting transaction data for Transaction Number: "
$wsConstants.Range("C10").Value = "Static part of logging message. Error retrieving Transaction Data."

# Row 11: LogMessage_Success (new position)
$wsConstants.Range("A11").Value = "LogMessage_Success"
$wsConstants.Range("B11").Value = "Transaction Successful."
$wsConstants.Range("C11").Value = "Static part of logging message. Processed Transaction succesful."

# Row 12: LogMessage_BusinessRuleException (new position, replacing old MailException row)
$wsConstants.Range("A12").Value = "LogMessage_BusinessRuleException"
$wsConstants.Range("B12").Value = "Business rule exception."
$wsConstants.Range("C12").Value = "Static part of logging message. Processed Transaction failed with business exception."

# Row 13: LogMessage_ApplicationException (new position)
$wsConstants.Range("A13").Value = "LogMessage_ApplicationException"
$wsConstants.Range("B13").Value = "System exception."
$wsConstants.Range("C13").Value = "Static part of logging message. Processed Transaction failed with application exception."

# Row 14: ExceptionMessage_ConsecutiveErrors (new)
$wsConstants.Range("A14").Value = "ExceptionMessage_ConsecutiveErrors"
$wsConstants.Range("B14").Value = "The maximum number of consecutive system exceptions was reached. "
$wsConstants.Range("C14").Value = "Error message in case MaxConsecutiveSystemExceptions number is reached."

# Row 15: stays blank separator

# Row 16: RetryNumberGetTransactionItem (new)
$wsConstants.Range("A16").Value = "RetryNumberGetTransactionItem"
$wsConstants.Range("B16").Value = 2
$wsConstants.Range("C16").Value = "The number of times Get Transaction Item activity is retried in case of an exception. Must be an integer >= 1."

# Row 17: RetryNumberSetTransactionStatus (new)
$wsConstants.Range("A17").Value = "RetryNumberSetTransactionStatus"
$wsConstants.Range("B17").Value = 2
$wsConstants.Range("C17").Value = "The number of times Set transaction status activity is retried in case of an exception. Must be an integer >= 1. "

# Row 18: stays blank separator

# Row 19: ShouldMarkJobAsFaulted (new)
$wsConstants.Range("A19").Value = "ShouldMarkJobAsFaulted"
$wsConstants.Range("B19").Value = $false
$wsConstants.Range("C19").Value = "Must be TRUE or FALSE. If the value is TRUE and an error occurs in Initialization state or the MaxConsecutiveSystemExceptions is reached, the job is marked as Faulted."

# ---------------------------------------------------------------------------
# Selections / active sheet: Constants becomes the active (selected) tab with
# B22 selected; Settings keeps A4 selected as its last cursor position.
# ---------------------------------------------------------------------------
$wsSettings.Activate()
$wsSettings.Range("A4").Select()

$wsConstants.Activate()
$wsConstants.Range("B22").Select()
